# Add tile movement limit to Disable Effects
# Adds two new KEY/ENGLISH rows to the "DisableEffects" worksheet:
#   TILE_MOVEMENT_MAX        -> "Max number of moving tiles at once:"
#   TILE_MOVEMENT_UNLIMITED  -> "Unlimited"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DisableEffects")

$ws.Cells.Item(9, 1).Value = "TILE_MOVEMENT_MAX"
$ws.Cells.Item(9, 2).Value = "Max number of moving tiles at once:"

$ws.Cells.Item(10, 1).Value = "TILE_MOVEMENT_UNLIMITED"
$ws.Cells.Item(10, 2).Value = "Unlimited"
